$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added New Mac-Address and Document Types: append a new data row (row 33)
# mirroring the existing regcntr_id / machine_id / lang_code / is_active /
# cr_by / cr_dtimes records.
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 10032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"

# Reposition the window/selection the way the workbook was left after the
# edit: scrolled down so row 19 is at the top, with E29 as the active cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("E29").Select()
